$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "image" column (E) values: drop the ".png" extension
$ws.Range("E2").Value = "project1"
$ws.Range("E3").Value = "project2"
$ws.Range("E4").Value = "project3"
$ws.Range("E5").Value = "project4"
$ws.Range("E6").Value = "project5"
$ws.Range("E7").Value = "project6"

# Move the selection/active cell to F8
$ws.Range("F8").Select()
